# Add "Poudre Paddlers" organization as a new row (row 4) in the "Data" sheet,
# shifting existing rows 4-6 down to 5-7, and fix up hyperlinks + view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remember the style used by the existing hyperlink/description cells (column E)
# so we can re-apply it after Hyperlinks.Add() mutates cell styles.
$linkStyle = $ws.Range("E2").Style

# ---------------------------------------------------------------------------
# 1. Insert a new blank row at position 4 (pushes old rows 4,5,6 -> 5,6,7)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new row 4 with the Poudre Paddlers data
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Poudre Paddlers"
$ws.Range("B4").Value = "Nonprofit"
$ws.Range("E4").Value = "https://www.poudrepaddlers.org/"
$ws.Range("C4").Value = "Canoe and kayak club"
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = -105.08152
$ws.Range("H4").Value = 40.583629999999999

$ws.Range("E4").Style = $linkStyle

# ---------------------------------------------------------------------------
# 3. Fix up the hyperlinks. Row-insert does not move the hyperlink anchors
#    in this engine, so the old hyperlinks are still sitting on E4/E5/E6
#    (now one row above where their data really is) and E3 needs to be
#    recreated last so relationship ids line up with the target layout.
# ---------------------------------------------------------------------------
function Remove-HyperlinksAtRow($ws, $rowNum) {
    $found = $true
    while ($found) {
        $found = $false
        foreach ($hh in $ws.Hyperlinks) {
            if ($hh.Range.Row -eq $rowNum) {
                $hh.Delete()
                $found = $true
                break
            }
        }
    }
}

Remove-HyperlinksAtRow $ws 4
Remove-HyperlinksAtRow $ws 5
Remove-HyperlinksAtRow $ws 6
Remove-HyperlinksAtRow $ws 3

$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.raftmw.com/")
$ws.Range("E7").Style = $linkStyle

$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.shoprma.com/")
$ws.Range("E6").Style = $linkStyle

$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.fcgov.com/parkplanning/poudre-river-park")
$ws.Range("E5").Style = $linkStyle

$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.americanwhitewater.org/")
$ws.Range("E3").Style = $linkStyle

# ---------------------------------------------------------------------------
# 4. Update the view/selection state for the Data sheet
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C5").Select()
